# Add a new "FinStats" dataset column (AB) to the Country-Year "Datasets and
# Years" tracker sheet: a header in AB1 and a checkmark for every data row
# (years 1990-2016 plus the two "varies"/"summarized" rows), mirroring the
# pattern already used by the neighbouring "UNCTAD Tariffs" column (AA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data -------------------------------------------------------------
# Header for the new last column.
$ws.Range("AB1").Value = "FinStats"

# Checkmark for every row that currently has one in column AA (rows 2-29).
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 28).Value = "$([char]0x2713)"
}

# --- formatting ---------------------------------------------------------
# Match the look of the rest of the header/data columns (Times New Roman
# 12pt, same as column AA) by copying AA's formatting onto the new AB cells.
$ws.Range("AA1:AA29").Copy()
$ws.Range("AB1:AB29").PasteSpecial(-4122)

# --- view state -----------------------------------------------------------
# Scroll the frozen pane over a couple of columns and move the active
# selection, same as the author did after adding the column.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 2
$ws.Range("Y12").Select()
